# Mark the e9ebefe9-615f-4844-b9d7-72d6ffb003e9 file as "Ready for handoff"
# (was "In Translation") and refresh its handoff timestamps, across the
# Overview summary sheet and each per-locale detail sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the e9ebefe9-615f-4844-b9d7-72d6ffb003e9.md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-14-18 08:14:39"

# --- zh-cn detail sheet: row 3 is the e9ebefe9-615f-4844-b9d7-72d6ffb003e9 row ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "2016-03-18 08:14:37"

# --- de-de detail sheet: row 3 is the e9ebefe9-615f-4844-b9d7-72d6ffb003e9 row ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "2016-03-18 08:14:39"
